$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "2026-01-30 10:06"
$ws.Range("B10").Value = 36
$ws.Range("C10").Value = 6
